# Refresh coin price / 1h-volume figures pulled from coinranking.com
# (scheduled GitHub Actions data sync)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.285.09"
$ws.Range("E2").Value = "  +4.74%  "

$ws.Range("D3").Value = "1.702.79"
$ws.Range("E3").Value = "  +4.14%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "'221.46"
$ws.Range("E5").Value = "  +2.82%  "

$ws.Range("E6").Value = "  +2.82%  "

$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "'29.80"
$ws.Range("E8").Value = "  +3.73%  "

$ws.Range("D9").Value = "'0.269"
$ws.Range("E9").Value = "  +3.07%  "

$ws.Range("E10").Value = "  +5.70%  "

$ws.Range("E11").Value = "  +1.21%  "

$ws.Range("D12").Value = "1.944.98"
$ws.Range("E12").Value = "  +4.01%  "

$ws.Range("D13").Value = "1.710.35"
$ws.Range("E13").Value = "  +4.51%  "

$ws.Range("E14").Value = "  +3.91%  "

$ws.Range("D15").Value = "'10.18"
$ws.Range("E15").Value = "  +8.30%  "

$ws.Range("D16").Value = "'4.18"
$ws.Range("E16").Value = "  +8.60%  "

$ws.Range("D17").Value = "31.267.09"
$ws.Range("E17").Value = "  +4.59%  "

$ws.Range("D18").Value = "'67.05"
$ws.Range("E18").Value = "  +3.85%  "

$ws.Range("D19").Value = "'250.69"
$ws.Range("E19").Value = "  +4.26%  "

$ws.Range("E20").Value = "  +3.30%  "

$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("E22").Value = "  +3.01%  "

$ws.Range("E23").Value = "  +3.46%  "

$ws.Range("E24").Value = "  -0.40%  "

$ws.Range("D25").Value = "'158.81"
$ws.Range("E25").Value = "  +0.96%  "

$ws.Range("D26").Value = "'16.03"
$ws.Range("E26").Value = "  +3.22%  "

$ws.Range("E27").Value = "  +3.13%  "

$ws.Range("D28").Value = "'6.80"
$ws.Range("E28").Value = "  +2.63%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("D30").Value = "'3.84"
$ws.Range("E30").Value = "  +13.51%  "

$ws.Range("E31").Value = "  +2.60%  "

$ws.Range("E32").Value = "  +3.86%  "

$ws.Range("D33").Value = "'3.41"
$ws.Range("E33").Value = "  +6.92%  "

$ws.Range("D34").Value = "1.509.95"
$ws.Range("E34").Value = "  +5.95%  "

$ws.Range("E35").Value = "  +3.08%  "

$ws.Range("E36").Value = "  +1.45%  "

$ws.Range("B37").Value = "Aave"
$ws.Range("C37").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D37").Value = "'83.22"
$ws.Range("E37").Value = "  +9.35%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "'0.614"
$ws.Range("E38").Value = "  +9.51%  "

$ws.Range("E39").Value = "  +4.76%  "

$ws.Range("E40").Value = "  -3.47%  "

$ws.Range("E41").Value = "  +0.41%  "

$ws.Range("D42").Value = "'2.05"
$ws.Range("E42").Value = "  +3.93%  "

$ws.Range("D43").Value = "'0.855"
$ws.Range("E43").Value = "  +2.73%  "

$ws.Range("D44").Value = "'0.0504"
$ws.Range("E44").Value = "  +0.77%  "

$ws.Range("E45").Value = "  +2.69%  "

$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").Value = "'52.31"
$ws.Range("E47").Value = "  +6.90%  "

$ws.Range("E48").Value = "  +3.42%  "

$ws.Range("D49").Value = "1.833.64"
$ws.Range("E49").Value = "  +3.19%  "

$ws.Range("E50").Value = "  +8.23%  "

$ws.Range("D51").Value = "'94.08"
$ws.Range("E51").Value = "  +1.31%  "

